# This script applies the "Updated cryptos list" data refresh to Sheet1.
# For each changed row, the Price (column D) and/or Volume(1h) (column E) cell
# is updated to its new value. The original cells store these as plain text
# (e.g. "29.056.19", "  -0.19%  "), so for any replacement text that Excel would
# otherwise auto-convert into a number, we force Text number-formatting before
# assigning the value and then restore the default "Normal" style afterwards so
# no visible formatting change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.056.19'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '1.835.01'
$ws.Range("E3").Value = '  -0.04%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9994'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '242.50'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.12%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.6143'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -3.35%  '

$ws.Range("E7").Value = '  -0.01%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.07476'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -1.08%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.2921'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("E10").Value = '  -1.12%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07685'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.55%  '

$ws.Range("D12").Value = '1.835.90'
$ws.Range("E12").Value = '  -0.01%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.994'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.31%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.6717'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.28%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '82.58'
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.000009146'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -4.34%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '5.909'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -3.09%  '

$ws.Range("D18").Value = '29.048.31'
$ws.Range("E18").Value = '  -0.34%  '

$ws.Range("D19").Value = '2.081.43'
$ws.Range("E19").Value = '  -0.28%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '231.51'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.32%  '

$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("E22").Value = '  +0.01%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '7.210'
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '158.74'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -1.16%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.1395'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -1.84%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '8.495'
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '17.80'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.98%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.495'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -0.64%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.149'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.41%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.124'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +1.14%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.05508'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.36%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.201'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.834'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -1.52%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.7370'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -1.47%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.142'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.17%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.659'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.25%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.772'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.54%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.01777'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.63%  '

$ws.Range("D40").Value = '1.211.76'
$ws.Range("E40").Value = '  -3.00%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '6.454'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -3.31%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.8905'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.73%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '101.88'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.24%  '

$ws.Range("D45").Value = '1.989.28'
$ws.Range("E45").Value = '  +0.14%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '65.41'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.08%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000121'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -3.20%  '

$ws.Range("E48").Value = '  -0.49%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.4062'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.38%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '9.088'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.62%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.05827'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.71%  '
